$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update G6: append completion note to the "Trể hẹn" status
$ws.Range("G6").Value = "Trể hẹn( Hoàn thành 20/10/2018)"

# Fill in G7 with the new completion status
$ws.Range("G7").Value = "Hoàn tthành(09/10/2018)"

# Widen column G to fit the new text (target stored width 31.7109375;
# the interop layer snaps ColumnWidth to 1/6-character steps, so
# 30.8333333333 is the closest input that lands on the nearest
# achievable stored width of 31.666666666666668)
$ws.Range("G1").ColumnWidth = 30.8333333333

# Move the active selection
$ws.Range("G22").Select()
